# Apply scraped cryptocurrency ranking updates to the sheet.
# Values that look numeric (e.g. "25.00", "0.7228") must be forced to
# literal text to match the source data (which stores prices as strings,
# including thousands-dot formatted figures like "29.366.34"), so a
# leading apostrophe is used to suppress Excel's automatic number
# detection, and the cell style is reset afterwards so no stray
# "quote prefix" / text-format styling is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.366.34'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.879.81'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'0.7228"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('D6').Value = "'243.04"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = "'0.08016"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('D9').Value = "'0.3144"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('D10').Value = "'25.00"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('D11').Value = "'0.08171"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.99%  '
$ws.Range('D12').Value = '1.884.38'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').Value = "'94.83"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.16%  '
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').Value = "'0.7136"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = "'6.426"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.88%  '
$ws.Range('D17').Value = "'0.000008516"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.26%  '
$ws.Range('D18').Value = '29.353.39'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').Value = "'244.48"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.125.82'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'0.9998"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = "'7.753"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = "'1.002"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = "'0.1605"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').Value = "'162.69"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = "'9.048"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'18.53"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = "'1.505"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = "'4.409"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = "'4.288"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = "'1.237"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.72%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.05365"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = "'1.937"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = "'0.7666"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.39%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = "'1.179"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = "'2.699"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'0.01873"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.264.86'
$ws.Range('E39').Value = '  +3.49%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = "'2.751"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'6.434"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = "'113.35"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.24%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'0.9081"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = "'74.59"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.15%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = "'0.00000000131"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.80%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = "'1.002"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '2.021.57'
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'0.5228"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = "'1.804"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'9.502"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').Value = "'0.4344"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.63%  '
